# Insert a new data row for "Femacal de La Calera" (Poroto verde) right
# before the current row 410, shifting the existing rows 410:506 down to
# 411:507 (dimension grows from A1:R506 to A1:R507).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 410 (and everything below it) down by one row.
$ws.Rows.Item(410).Insert()

# Populate the newly-inserted row 410 with its values.
$ws.Cells.Item(410, 1).Value = 3
$ws.Cells.Item(410, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(410, 3).Value = "Coquimbo"
$ws.Cells.Item(410, 4).Value = 44943
$ws.Cells.Item(410, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(410, 5).Value = 5
$ws.Cells.Item(410, 6).Value = 100112031
$ws.Cells.Item(410, 7).Value = "Poroto verde"
$ws.Cells.Item(410, 8).Value = "Magnum"
$ws.Cells.Item(410, 9).Value = "Primera"
$ws.Cells.Item(410, 10).Value = 78
$ws.Cells.Item(410, 11).Value = 25000
$ws.Cells.Item(410, 12).Value = 26000
$ws.Cells.Item(410, 13).Value = 25487
$ws.Cells.Item(410, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(410, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(410, 16).Value = 1019
$ws.Cells.Item(410, 17).Value = 25
$ws.Cells.Item(410, 18).Value = "Hortaliza"
